$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (CPF), shifting CPF -> J, CEP -> K
$ws.Columns.Item(9).Insert()

# Update header row
$ws.Range("I1").Value = "Celular"

# The CEP column (K) holds numeric-looking codes (some with leading zeros),
# so force it to be stored as text to avoid Excel auto-converting to a number.
$ws.Range("K2:K11").NumberFormat = "@"

# Row 2
$ws.Range("B2").Value = "Rafael Hugo Silveira Ribeira Ferreira"
$ws.Range("C2").Value = 4
$ws.Range("F2").Value = "Cecília Juliana Ferreira Silveira"
$ws.Range("H2").Value = "Fora da força de trabalho"
$ws.Range("I2").Value = "(48) 97605-3870"
$ws.Range("J2").Value = "887.612.030-09"
$ws.Range("K2").Value = "921188030"

# Row 3
$ws.Range("B3").Value = "Carlos Fernando Ribeiro Nunes"
$ws.Range("C3").Value = 43
$ws.Range("D3").Value = "Masculino"
$ws.Range("E3").Value = "Pardo"
$ws.Range("F3").Value = "Priscila Nunes Fernando"
$ws.Range("G3").Value = "Fundamental incompleto"
$ws.Range("I3").Value = "(85) 95121-8739"
$ws.Range("J3").Value = "710.090.900-70"
$ws.Range("K3").Value = "202427825"

# Row 4
$ws.Range("B4").Value = "Lorena Dias"
$ws.Range("C4").Value = 22
$ws.Range("E4").Value = "Branca"
$ws.Range("F4").Value = "Flávia Monteiro Dias"
$ws.Range("I4").Value = "(73) 96202-0181"
$ws.Range("J4").Value = "228.714.090-56"
$ws.Range("K4").Value = "103240838"

# Row 5
$ws.Range("B5").Value = "Hugo Gustavo Ramos Campos"
$ws.Range("C5").Value = 18
$ws.Range("E5").Value = "Branco"
$ws.Range("F5").Value = "Sueli Cecília Pereira Campos Gustavo Ramos"
$ws.Range("H5").Value = "Empregado: Setor privado (CLT)"
$ws.Range("I5").Value = "(92) 99506-7640"
$ws.Range("J5").Value = "368.774.670-43"
$ws.Range("K5").Value = "251849819"

# Row 6
$ws.Range("B6").Value = "Joaquim Nogueira"
$ws.Range("C6").Value = 39
$ws.Range("D6").Value = "Masculino"
$ws.Range("E6").Value = "Preto"
$ws.Range("F6").Value = "Carolina Ribeira Nogueira"
$ws.Range("G6").Value = "Médio incompleto"
$ws.Range("H6").Value = "Empregado: Setor privado (CLT)"
$ws.Range("I6").Value = "(87) 93734-5723"
$ws.Range("J6").Value = "991.471.980-52"
$ws.Range("K6").Value = "002896322"

# Row 7
$ws.Range("B7").Value = "Armando Dias Melo Alves Monteiro"
$ws.Range("C7").Value = 11
$ws.Range("D7").Value = "Masculino"
$ws.Range("E7").Value = "Pardo"
$ws.Range("F7").Value = "Bianca Dias Melo"
$ws.Range("H7").Value = "Fora da força de trabalho"
$ws.Range("I7").Value = "(46) 99951-2572"
$ws.Range("J7").Value = "914.311.810-05"
$ws.Range("K7").Value = "969481242"

# Row 8
$ws.Range("B8").Value = "Luiza Sandra Alves Andrade Azevedo"
$ws.Range("C8").Value = 49
$ws.Range("F8").Value = "Bianca Vanessa Azevedo Sandra Andrade"
$ws.Range("G8").Value = "Médio completo"
$ws.Range("I8").Value = "(53) 96393-7539"
$ws.Range("J8").Value = "029.707.070-39"
$ws.Range("K8").Value = "687937125"

# Row 9
$ws.Range("B9").Value = "Bruno Sousa"
$ws.Range("C9").Value = 11
$ws.Range("D9").Value = "Masculino"
$ws.Range("E9").Value = "Outro"
$ws.Range("F9").Value = "Luiza Machado Borges Pires Sousa"
$ws.Range("G9").Value = "Em idade escolar."
$ws.Range("H9").Value = "Fora da força de trabalho"
$ws.Range("I9").Value = "(48) 98175-6393"
$ws.Range("J9").Value = "187.167.750-58"
$ws.Range("K9").Value = "227263497"

# Row 10
$ws.Range("B10").Value = "Juliana Cavalcanti Pinto Martins"
$ws.Range("C10").Value = 10
$ws.Range("D10").Value = "Feminino"
$ws.Range("E10").Value = "Branca"
$ws.Range("F10").Value = "Nathalia Elaine Pereira Pinto Cavalcanti Martins"
$ws.Range("I10").Value = "(44) 98296-7618"
$ws.Range("J10").Value = "319.678.350-20"
$ws.Range("K10").Value = "179029847"

# Row 11
$ws.Range("B11").Value = "Natália Silveira Cavalcanti"
$ws.Range("C11").Value = 54
$ws.Range("E11").Value = "Parda"
$ws.Range("F11").Value = "Laura Cavalcanti Silveira"
$ws.Range("G11").Value = "Fundamental incompleto"
$ws.Range("H11").Value = "Autonomo: Sem CNPJ"
$ws.Range("I11").Value = "(85) 90785-2366"
$ws.Range("J11").Value = "471.400.420-47"
$ws.Range("K11").Value = "267793507"
